# Generate Report for Handback
# Update the status / handback datetime / error detail for the
# "c6e1dc52-95da-4f96-84be-ceed4a84112a.md" row now that the handback
# has been generated and is in sync with en-US.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-19 08:46:26"
$zhcn.Range("P3").Value = ""

# --- de-de sheet --------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-19 08:46:32"
$dede.Range("P3").Value = ""
